$d = $word.ActiveDocument

# 1) Merge the two "caption" / " reference" runs (same PSRefCaption style) into a
#    single run reading "caption reference" (simple re-type over the found text,
#    which Word collapses into one run).
$d.Content.Find.Execute("caption reference", $false, $false, $false, $false, $false, $true, 1, $false, "caption reference", 2) | Out-Null

# 2) Sections 2, 3 and 4 currently inherit their header from the previous section
#    (LinkToPrevious = True), which is why no header shows up for them on export.
#    Give each of them (and section 4's footer) their own explicit header/footer,
#    matching the content that would otherwise have been inherited.
for ($i = 2; $i -le $d.Sections.Count; $i++) {
    $section = $d.Sections.Item($i)

    $header = $section.Headers.Item(1)
    if ($header.LinkToPrevious) {
        $header.LinkToPrevious = $false
    }
}

$lastSection = $d.Sections.Item($d.Sections.Count)
$lastFooter = $lastSection.Footers.Item(1)
if ($lastFooter.LinkToPrevious) {
    $lastFooter.LinkToPrevious = $false
}

Write-Output "done"
